$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header cells to the new terminology
$ws.Range("I1").Value = "Market 1 Fiat Spot Price"
$ws.Range("J1").Value = "Market 2 Fiat Spot Price"
$ws.Range("K1").Value = "Fee Asset Fiat Spot Price"
$ws.Range("H1").Value = "Fee Asset"

# Restore the view: scroll to show row 1, and select H2
$ws.Activate()
$ws.Range("A1").Select()
$ws.Range("H2").Select()
